$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels: *_old -> *_FV2404, *_new -> *_FV2410 ---------
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- 2. Freeze the header row (pane split below row 1) ---------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into an Excel Table -----------------------------
# Stash the header row's existing formatting on a scratch row, reset the
# header to the default style so ListObjects.Add doesn't snapshot the bold
# header look into a new dxf/table style, then restore the formatting once
# the table exists.
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A200:U200")
$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

$headerRange.Style = "Normal"

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U78"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = $null

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()
$scratchRange.Clear()

Write-Host "Done"
